$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(6,  "Nabhya",        "nabhyakhoria@gmail.com", "ABCD",                 "BABACSZCascs"),
    @(7,  "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(8,  "nabhyakhoria@gmail.com", "nabhyakhoria@gmail.com", "sadsad",      "wdqd"),
    @(9,  "nabhyakhoria@gmail.com", "nabhyakhoria@gmail.com", "sadsad",      "wdqd"),
    @(10, "nabhyakhoria@gmail.com", "nabhyakhoria@gmail.com", "sadsad",      "wdqd"),
    @(11, "nabhyakhoria@gmail.com", "nabhyakhoria@gmail.com", "sadsad",      "wdqd"),
    @(12, "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(13, "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(14, "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(15, "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(16, "Nabhya",        "nabhyakhoria@gmail.com", "Room fan not working", "Room SDS 440.`nCheck asap"),
    @(17, "Nabhya Khoria", "nabhyakhoria@gmail.com", "Nice",                 "HAt"),
    @(18, "Shivam Singh",  "shivam@gmail.com",       "VP",                   "New VP ")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
